$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as plain text in the workbook (e.g.
# "0.998", "320.99", "2.634.17"). Some of the updated prices below are
# round numbers like "1.00", "322.00" or "125.21" which Excel would silently
# reinterpret as numeric values (dropping the trailing zeros / changing the
# representation) unless the cell is explicitly formatted as Text first.
$textCells = "D4","D5","D6","D7","D8","D9","D10","D11","D14","D17","D19","D20","D23","D24","D25","D26","D28","D30","D32","D34","D36","D37","D40","D41","D43","D44","D46","D47","D50","D51"
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '48.936.08'
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").Value = '2.638.18'
$ws.Range("E3").Value = '  +5.44%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '110.47'
$ws.Range("E5").Value = '  +4.01%  '
$ws.Range("D6").Value = '322.00'
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("D7").Value = '0.519'
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").Value = '0.540'
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").Value = '39.50'
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").Value = '19.85'
$ws.Range("E11").Value = '  -1.30%  '
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").Value = '7.21'
$ws.Range("E14").Value = '  +1.55%  '
$ws.Range("D15").Value = '3.048.87'
$ws.Range("E15").Value = '  +5.47%  '
$ws.Range("D16").Value = '2.638.86'
$ws.Range("E16").Value = '  +6.66%  '
$ws.Range("D17").Value = '0.858'
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").Value = '49.009.09'
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("D20").Value = '6.66'
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("D22").Value = '0.0₃0941'
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("D23").Value = '269.79'
$ws.Range("E23").Value = '  -3.58%  '
$ws.Range("D24").Value = '70.02'
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("D25").Value = '2.55'
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("D26").Value = '26.22'
$ws.Range("E26").Value = '  +1.69%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = '10.06'
$ws.Range("E28").Value = '  +4.12%  '
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("D30").Value = '35.11'
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("E31").Value = '  -1.75%  '
$ws.Range("D32").Value = '49.32'
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("E33").Value = '  +2.86%  '
$ws.Range("D34").Value = '19.19'
$ws.Range("E34").Value = '  -1.20%  '
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("D36").Value = '0.0793'
$ws.Range("E36").Value = '  +2.79%  '
$ws.Range("D37").Value = '4.92'
$ws.Range("E37").Value = '  +8.63%  '
$ws.Range("E38").Value = '  +3.67%  '
$ws.Range("E39").Value = '  +8.56%  '
$ws.Range("D40").Value = '125.21'
$ws.Range("D41").Value = '22.51'
$ws.Range("E41").Value = '  +4.12%  '
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("D43").Value = '2.18'
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").Value = '0.0312'
$ws.Range("E44").Value = '  +3.82%  '
$ws.Range("D45").Value = '2.083.19'
$ws.Range("E45").Value = '  +4.63%  '
$ws.Range("D46").Value = '3.21'
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("D47").Value = '2.11'
$ws.Range("E47").Value = '  +7.37%  '
$ws.Range("E48").Value = '  +4.71%  '
$ws.Range("D49").Value = '2.882.72'
$ws.Range("E49").Value = '  +5.12%  '
$ws.Range("D50").Value = '8.90'
$ws.Range("E50").Value = '  -0.92%  '
$ws.Range("D51").Value = '59.05'
$ws.Range("E51").Value = '  +4.95%  '
